$d = $word.ActiveDocument
$d.Content.Find.Execute("Documents_for_", $true, $false, $false, $false, $false,
                         $true, 1, $false, "documents_for_", 2)
